# Change the table style (table design/theme) applied to the table on
# Slide 5 from the default "Table_0" style ({A2A8FE15-AAEA-4D18-B53B-175FA3B5CFA0})
# to the built-in table style {8FE565B8-BC63-44E3-9304-7200336F656F}.

$p = $ppt.ActivePresentation

$newStyleId = "{8FE565B8-BC63-44E3-9304-7200336F656F}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
